$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3  (Id 111943877)
$ws.Range("A3").Value = 111943877
$ws.Range("B3").Value = 90666
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 4364
$ws.Range("F3").Value = "Dropptaggsvamp"
$ws.Range("G3").Value = "Hydnellum ferrugineum"
$ws.Range("H3").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q3").Value = 600476
$ws.Range("R3").Value = 7221499
$ws.Range("AX3").Value = "Simon Mattsson, Maja Östlund"

# Row 4  (Id 111943803)
$ws.Range("A4").Value = 111943803
$ws.Range("B4").Value = 89405
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 1202
$ws.Range("F4").Value = "Ullticka"
$ws.Range("G4").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H4").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q4").Value = 600424
$ws.Range("R4").Value = 7221684
$ws.Range("AX4").Value = "Simon Mattsson, Maja Östlund"

# Row 5  (Id 111943882)
$ws.Range("A5").Value = 111943882
$ws.Range("B5").Value = 90666
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 4364
$ws.Range("F5").Value = "Dropptaggsvamp"
$ws.Range("G5").Value = "Hydnellum ferrugineum"
$ws.Range("H5").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q5").Value = 600419
$ws.Range("R5").Value = 7221432
$ws.Range("AX5").Value = "Simon Mattsson, Maja Östlund"

# Row 6  (Id 111943841)
$ws.Range("A6").Value = 111943841
$ws.Range("B6").Value = 77268
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 228912
$ws.Range("F6").Value = "Mörk kolflarnlav"
$ws.Range("G6").Value = "Carbonicola myrmecina"
$ws.Range("H6").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q6").Value = 600367
$ws.Range("R6").Value = 7221297
$ws.Range("AX6").Value = "Simon Mattsson, Maja Östlund"

# Row 7  (Id 111943881)
$ws.Range("A7").Value = 111943881
$ws.Range("B7").Value = 90666
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 4364
$ws.Range("F7").Value = "Dropptaggsvamp"
$ws.Range("G7").Value = "Hydnellum ferrugineum"
$ws.Range("H7").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q7").Value = 600419
$ws.Range("R7").Value = 7221630
$ws.Range("AX7").Value = "Simon Mattsson, Maja Östlund"

# Row 8  (Id 111943944)
$ws.Range("A8").Value = 111943944
$ws.Range("B8").Value = 90666
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 4364
$ws.Range("F8").Value = "Dropptaggsvamp"
$ws.Range("G8").Value = "Hydnellum ferrugineum"
$ws.Range("H8").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q8").Value = 600345
$ws.Range("R8").Value = 7221526
$ws.Range("AX8").Value = "Maja Östlund, Simon Mattsson"

# Row 9  (Id 111943887)
$ws.Range("A9").Value = 111943887
$ws.Range("B9").Value = 90660
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 4362
$ws.Range("F9").Value = "Blå taggsvamp"
$ws.Range("G9").Value = "Hydnellum caeruleum"
$ws.Range("H9").Value = "(Hornem.) P.Karst."
$ws.Range("Q9").Value = 600485
$ws.Range("R9").Value = 7221470
$ws.Range("AX9").Value = "Simon Mattsson, Maja Östlund"

# Row 10  (Id 111943940)
$ws.Range("A10").Value = 111943940
$ws.Range("B10").Value = 77515
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("Q10").Value = 600237
$ws.Range("R10").Value = 7221447
$ws.Range("AX10").Value = "Maja Östlund, Simon Mattsson"

# Row 11  (Id 111943907)
$ws.Range("A11").Value = 111943907
$ws.Range("B11").Value = 90682
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 2059
$ws.Range("F11").Value = "Skrovlig taggsvamp"
$ws.Range("G11").Value = "Hydnellum scabrosum"
$ws.Range("H11").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q11").Value = 600409
$ws.Range("R11").Value = 7221648
$ws.Range("AX11").Value = "Maja Östlund, Simon Mattsson"

# Row 12  (Id 111943879)
$ws.Range("A12").Value = 111943879
$ws.Range("B12").Value = 90666
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 4364
$ws.Range("F12").Value = "Dropptaggsvamp"
$ws.Range("G12").Value = "Hydnellum ferrugineum"
$ws.Range("H12").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q12").Value = 600452
$ws.Range("R12").Value = 7221545
$ws.Range("AX12").Value = "Simon Mattsson, Maja Östlund"

# Row 13  (Id 111943815)
$ws.Range("A13").Value = 111943815
$ws.Range("B13").Value = 90682
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 2059
$ws.Range("F13").Value = "Skrovlig taggsvamp"
$ws.Range("G13").Value = "Hydnellum scabrosum"
$ws.Range("H13").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q13").Value = 600430
$ws.Range("R13").Value = 7221629
$ws.Range("AX13").Value = "Simon Mattsson, Maja Östlund"

# Row 14  (Id 111943880)
$ws.Range("A14").Value = 111943880
$ws.Range("B14").Value = 90666
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 4364
$ws.Range("F14").Value = "Dropptaggsvamp"
$ws.Range("G14").Value = "Hydnellum ferrugineum"
$ws.Range("H14").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q14").Value = 600437
$ws.Range("R14").Value = 7221453
$ws.Range("AX14").Value = "Simon Mattsson, Maja Östlund"

# Row 15  (Id 111943947)
$ws.Range("A15").Value = 111943947
$ws.Range("B15").Value = 85715
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 510
$ws.Range("F15").Value = "Doftskinn"
$ws.Range("G15").Value = "Cystostereum murrayi"
$ws.Range("H15").Value = "(Berk. & M.A. Curtis.) Pouzar"
$ws.Range("Q15").Value = 600352
$ws.Range("R15").Value = 7221402
$ws.Range("AX15").Value = "Maja Östlund, Simon Mattsson"

# Row 16  (Id 111943883)
$ws.Range("A16").Value = 111943883
$ws.Range("B16").Value = 90666
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 4364
$ws.Range("F16").Value = "Dropptaggsvamp"
$ws.Range("G16").Value = "Hydnellum ferrugineum"
$ws.Range("H16").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q16").Value = 600311
$ws.Range("R16").Value = 7221358
$ws.Range("AX16").Value = "Simon Mattsson, Maja Östlund"

# Row 17  (Id 111943814)
$ws.Range("A17").Value = 111943814
$ws.Range("B17").Value = 90682
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 2059
$ws.Range("F17").Value = "Skrovlig taggsvamp"
$ws.Range("G17").Value = "Hydnellum scabrosum"
$ws.Range("H17").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q17").Value = 600437
$ws.Range("R17").Value = 7221630
$ws.Range("AX17").Value = "Simon Mattsson, Maja Östlund"

# Row 18  (Id 111943816)
$ws.Range("A18").Value = 111943816
$ws.Range("B18").Value = 90682
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 2059
$ws.Range("F18").Value = "Skrovlig taggsvamp"
$ws.Range("G18").Value = "Hydnellum scabrosum"
$ws.Range("H18").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q18").Value = 600428
$ws.Range("R18").Value = 7221623
$ws.Range("AX18").Value = "Simon Mattsson, Maja Östlund"

# Remove the Starttid / Sluttid (time) columns for data rows 3-18
$ws.Range("Z3:Z18").ClearContents()
$ws.Range("AB3:AB18").ClearContents()
